# The automatic re-export reordered three observation records (rows 4-6 on
# sheet "Artfynd"). The row payloads rotate: old row 4 -> new row 5,
# old row 5 -> new row 6, old row 6 -> new row 4. A handful of columns
# (A, P, Q, R, S, AW, AX) differ between the three records, plus a set of
# blank placeholder cells (J, L, N, AF) that old rows 5/6 carried but old
# row 4 didn't - after the rotation, new row 4 needs them and new row 5
# loses them (new row 6 keeps them, since it inherits old row 5's shape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix up the blank placeholder cells first (before row 6's own values change) ---
# New row 4 = old row 6's shape -> bring in its blank J/L/N/AF placeholders.
$ws.Cells.Item(6, 10).Copy($ws.Cells.Item(4, 10))  # J4
$ws.Cells.Item(6, 12).Copy($ws.Cells.Item(4, 12))  # L4
$ws.Cells.Item(6, 14).Copy($ws.Cells.Item(4, 14))  # N4
$ws.Cells.Item(6, 32).Copy($ws.Cells.Item(4, 32))  # AF4

# New row 5 = old row 4's shape -> it never had J/L/N/AF, so drop them.
$ws.Cells.Item(5, 10).ClearContents()  # J5
$ws.Cells.Item(5, 12).ClearContents()  # L5
$ws.Cells.Item(5, 14).ClearContents()  # N5
$ws.Cells.Item(5, 32).ClearContents()  # AF5

# --- rotate the per-record values ---
# Row 4 becomes what used to be row 6 ("Albinvägen2...")
$ws.Range("A4").Value = 111454300
$ws.Range("P4").Value = "Albinvägen2, Svartnäset, Hackås, Jmt"
$ws.Range("Q4").Value = 492448.9318965223
$ws.Range("R4").Value = 6948282.559996245
$ws.Range("S4").Value = 10
$ws.Range("AW4").Value = "Jan Magnesved"
$ws.Range("AX4").Value = "Jan Magnesved, Anders Wännström "

# Row 5 becomes what used to be row 4 ("Siljebodarna...")
$ws.Range("A5").Value = 111454959
$ws.Range("P5").Value = "Siljebodarna, Jmt"
$ws.Range("Q5").Value = 492425.096130528
$ws.Range("R5").Value = 6948324.435442663
$ws.Range("S5").Value = 15
$ws.Range("AW5").Value = "Monica Magnesved"
$ws.Range("AX5").Value = "Monica Magnesved"

# Row 6 becomes what used to be row 5 ("Albinvägen3...")
$ws.Range("A6").Value = 111454321
$ws.Range("P6").Value = "Albinvägen3, Svartnäset, Hackås, Jmt"
$ws.Range("Q6").Value = 492408.8681431987
$ws.Range("R6").Value = 6948272.081593725
$ws.Range("S6").Value = 10
$ws.Range("AW6").Value = "Jan Magnesved"
$ws.Range("AX6").Value = "Jan Magnesved, Anders Wännström "
